# Scheduled-runner market-data refresh: update the market/price-derived
# columns (H:N - currentAveragePrice*, LevePrice*, LeveProfit*) for the
# leves whose market data changed, across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# sheets. Column A:G (leve name/item/level/exp/gil/amount/itemId) are left
# untouched since they don't change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2289.0908
$ws.Range("I51").Value = 1260
$ws.Range("J51").Value = 3146.6667
$ws.Range("K51").Value = 1260
$ws.Range("L51").Value = 3146.6667
$ws.Range("M51").Value = -776
$ws.Range("N51").Value = -4114.6667

$ws.Range("H138").Value = 15509571
$ws.Range("I138").Value = 41671372
$ws.Range("J138").Value = 6280
$ws.Range("K138").Value = 125014116
$ws.Range("L138").Value = 18840
$ws.Range("M138").Value = -125008976
$ws.Range("N138").Value = -29120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25827.684
$ws.Range("I32").Value = 25827.684
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 25827.684
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -25540.684
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 1338324.5
$ws.Range("I45").Value = 1977644.6
$ws.Range("K45").Value = 1977644.6
$ws.Range("M45").Value = -1977267.6

$ws.Range("H61").Value = 2042.3422
$ws.Range("I61").Value = 2042.3422
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2042.3422
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1830.3422
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 6177.857
$ws.Range("I74").Value = 1345.7222
$ws.Range("J74").Value = 35170.668
$ws.Range("K74").Value = 1345.7222
$ws.Range("L74").Value = 35170.668
$ws.Range("M74").Value = -471.7221999999999
$ws.Range("N74").Value = -36918.668

$ws.Range("H77").Value = 6177.857
$ws.Range("I77").Value = 1345.7222
$ws.Range("J77").Value = 35170.668
$ws.Range("K77").Value = 6728.611
$ws.Range("L77").Value = 175853.34
$ws.Range("M77").Value = -2360.611
$ws.Range("N77").Value = -184589.34

$ws.Range("H102").Value = 1773.091
$ws.Range("I102").Value = 1500.8889
$ws.Range("J102").Value = 2998
$ws.Range("K102").Value = 1500.8889
$ws.Range("L102").Value = 2998
$ws.Range("M102").Value = 121.1111000000001
$ws.Range("N102").Value = -6242

$ws.Range("H136").Value = 2042.3422
$ws.Range("I136").Value = 2042.3422
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6127.0266
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3577.0266
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 59998
$ws.Range("J137").Value = 59998
$ws.Range("L137").Value = 59998
$ws.Range("N137").Value = -70198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 894.21875
$ws.Range("I94").Value = 772.88464
$ws.Range("J94").Value = 1420
$ws.Range("K94").Value = 772.88464
$ws.Range("L94").Value = 1420
$ws.Range("M94").Value = -321.88464
$ws.Range("N94").Value = -2322

$ws.Range("H118").Value = 7855.8823
$ws.Range("J118").Value = 7855.8823
$ws.Range("L118").Value = 7855.8823
$ws.Range("N118").Value = -11169.8823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5636991
$ws.Range("I6").Value = 5200690
$ws.Range("J6").Value = 10000000
$ws.Range("K6").Value = 5200690
$ws.Range("L6").Value = 10000000
$ws.Range("M6").Value = -5200577
$ws.Range("N6").Value = -10000226

$ws.Range("H7").Value = 139.2
$ws.Range("I7").Value = 45.23077
$ws.Range("J7").Value = 750
$ws.Range("K7").Value = 45.23077
$ws.Range("L7").Value = 750
$ws.Range("M7").Value = 67.76922999999999
$ws.Range("N7").Value = -976

$ws.Range("H41").Value = 30016.25
$ws.Range("J41").Value = 30016.25
$ws.Range("L41").Value = 30016.25
$ws.Range("N41").Value = -30872.25

$ws.Range("H50").Value = 9174.5
$ws.Range("J50").Value = 9174.5
$ws.Range("L50").Value = 9174.5
$ws.Range("N50").Value = -10424.5

$ws.Range("H51").Value = 9339.6
$ws.Range("J51").Value = 9339.6
$ws.Range("L51").Value = 9339.6
$ws.Range("N51").Value = -10811.6

$ws.Range("H59").Value = 15731.75
$ws.Range("J59").Value = 15731.75
$ws.Range("L59").Value = 15731.75
$ws.Range("N59").Value = -18021.75

$ws.Range("H60").Value = 7637.75
$ws.Range("I60").Value = 6666.6665
$ws.Range("J60").Value = 8220.4
$ws.Range("K60").Value = 6666.6665
$ws.Range("L60").Value = 8220.4
$ws.Range("M60").Value = -6155.6665
$ws.Range("N60").Value = -9242.4

$ws.Range("H61").Value = 9339.6
$ws.Range("J61").Value = 9339.6
$ws.Range("L61").Value = 9339.6
$ws.Range("N61").Value = -10035.6

$ws.Range("H68").Value = 17374.875
$ws.Range("J68").Value = 17374.875
$ws.Range("L68").Value = 17374.875
$ws.Range("N68").Value = -18872.875

$ws.Range("H71").Value = 17374.875
$ws.Range("J71").Value = 17374.875
$ws.Range("L71").Value = 52124.625
$ws.Range("N71").Value = -59612.625

$ws.Range("H74").Value = 13999.75
$ws.Range("J74").Value = 13999.75
$ws.Range("L74").Value = 13999.75
$ws.Range("N74").Value = -15747.75

$ws.Range("H77").Value = 13999.75
$ws.Range("J77").Value = 13999.75
$ws.Range("L77").Value = 41999.25
$ws.Range("N77").Value = -50735.25

$ws.Range("H105").Value = 1716
$ws.Range("I105").Value = 1734
$ws.Range("K105").Value = 1734
$ws.Range("M105").Value = 13

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 644.5454999999999
$ws.Range("I34").Value = 136.25
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 408.75
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -324.75
$ws.Range("N34").Value = -6168

$ws.Range("H39").Value = 496.42856

$ws.Range("H55").Value = 425.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 715.8570999999999
$ws.Range("I46").Value = 466.66666
$ws.Range("J46").Value = 902.75
$ws.Range("K46").Value = 466.66666
$ws.Range("L46").Value = 902.75
$ws.Range("M46").Value = -278.66666
$ws.Range("N46").Value = -1278.75

$ws.Range("H55").Value = 342
$ws.Range("I55").Value = 196.92308
$ws.Range("J55").Value = 513.4545000000001
$ws.Range("K55").Value = 196.92308
$ws.Range("L55").Value = 513.4545000000001
$ws.Range("M55").Value = -23.92308
$ws.Range("N55").Value = -859.4545000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3300
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 4966.6665
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 9933.333000000001
$ws.Range("M81").Value = -539
$ws.Range("N81").Value = -12055.333

$ws.Range("H84").Value = 3300
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 4966.6665
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 49666.665
$ws.Range("M84").Value = -2696
$ws.Range("N84").Value = -60274.665

$ws.Range("H136").Value = 1361.94
$ws.Range("I136").Value = 963.05
$ws.Range("J136").Value = 2957.5
$ws.Range("K136").Value = 2889.15
$ws.Range("L136").Value = 8872.5
$ws.Range("M136").Value = -339.1499999999996
$ws.Range("N136").Value = -13972.5
